$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume/1h (E) columns for updated crypto data ---

# Row 2
$ws.Range("D2").Value = "30.117.94"
$ws.Range("E2").Value = "  +0.54%  "

# Row 3
$ws.Range("D3").Value = "1.884.61"
$ws.Range("E3").Value = "  +0.24%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.41"
$ws.Range("E5").Value = "  -2.27%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9987"
$ws.Range("E6").Value = "  -0.31%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4977"
$ws.Range("E7").Value = "  +0.20%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2934"
$ws.Range("E8").Value = "  +2.86%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06613"
$ws.Range("E9").Value = "  +1.03%  "

# Row 10
$ws.Range("D10").Value = "1.880.68"
$ws.Range("E10").Value = "  -0.06%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.87"
$ws.Range("E11").Value = "  -1.37%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07177"
$ws.Range("E12").Value = "  -0.22%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6639"
$ws.Range("E13").Value = "  +0.24%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "85.80"
$ws.Range("E14").Value = "  +0.78%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.858"
$ws.Range("E15").Value = "  +1.57%  "

# Row 16
$ws.Range("D16").Value = "30.105.27"
$ws.Range("E16").Value = "  +0.51%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007923"
$ws.Range("E17").Value = "  +5.76%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9989"
$ws.Range("E18").Value = "  -0.27%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("E19").Value = "  -0.50%  "

# Row 20
$ws.Range("D20").Value = "2.121.79"
$ws.Range("E20").Value = "  -0.15%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9980"
$ws.Range("E21").Value = "  -0.28%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.757"
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.606"
$ws.Range("E23").Value = "  +1.41%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.127"
$ws.Range("E24").Value = "  +1.48%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.88"
$ws.Range("E25").Value = "  +5.28%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "134.63"
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.79"
$ws.Range("E27").Value = "  +0.58%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.915"
$ws.Range("E28").Value = "  -2.33%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.381"
$ws.Range("E29").Value = "  -1.25%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.167"
$ws.Range("E30").Value = "  -0.05%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08690"
$ws.Range("E31").Value = "  +0.99%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.947"
$ws.Range("E32").Value = "  +1.95%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05004"
$ws.Range("E33").Value = "  -1.56%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7071"
$ws.Range("E34").Value = "  +3.36%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.107"
$ws.Range("E35").Value = "  -2.29%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.649"
$ws.Range("E36").Value = "  -2.50%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.701"
$ws.Range("E37").Value = "  -1.80%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9357"
$ws.Range("E39").Value = "  -2.51%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01650"
$ws.Range("E40").Value = "  +1.34%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.962"
$ws.Range("E41").Value = "  -1.86%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9990"
$ws.Range("E42").Value = "  -0.32%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.51"
$ws.Range("E44").Value = "  -1.24%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.493"
$ws.Range("E45").Value = "  +0.99%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05695"
$ws.Range("E47").Value = "  +1.16%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "32.44"
$ws.Range("E48").Value = "  +0.34%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.257"
$ws.Range("E49").Value = "  +0.26%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.344"
$ws.Range("E50").Value = "  +0.81%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3717"
$ws.Range("E51").Value = "  +0.32%  "

# Rows where only the Volume(1h) column changed
# Row 38
$ws.Range("E38").Value = "  -4.50%  "

# Row 43
$ws.Range("E43").Value = "  -0.26%  "

# Row 46
$ws.Range("E46").Value = "  +0.37%  "

